$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.205.48"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.596.03"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("D5").Value = "'212.09"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'18.90"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'0.0854"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "1.580.72"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "'0.505"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "'63.60"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "26.220.98"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'226.73"
$ws.Range("E18").Value = "  +5.77%  "
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'8.90"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "'145.65"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "'15.34"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "'0.0492"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "1.443.73"
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "'0.564"
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "'0.818"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D42").Value = "'2.17"
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("D43").Value = "'0.930"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").Value = "1.734.69"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'0.757"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").Value = "'60.37"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "'87.57"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("D48").Value = "'1.47"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'0.0946"
$ws.Range("E51").Value = "  -3.14%  "
